$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(77, 1).Value = 27
$ws.Cells.Item(77, 2).Value = 329
$ws.Cells.Item(77, 3).Value = 1616
